# Regenerate the localization-status report for the new handoff commit
# (old source commit guid cb648a3a-7038-401f-89ec-1d8a7ceec1f0 ->
#  new source commit guid 8bfdab98-99e5-42de-8fb8-e88c4881a2b4) and record
# that handback has not happened yet for this new handoff.

$wb = $excel.ActiveWorkbook

$oldGuid = "cb648a3a-7038-401f-89ec-1d8a7ceec1f0"
$newGuid = "8bfdab98-99e5-42de-8fb8-e88c4881a2b4"
$oldHash = "6e7625768682252bf0fc1fd4865b23bb2d033112"
$newHash = "37b344794930aefebefce6e82a4e80297a910ea8"

$newHoGenerateDate = "2016-08-19 17:03:34"
$zhHandoffDate     = "2016-08-19 17:03:30"
$noHandbackDate    = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Hyperlinks.Delete()
$overviewHlTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/$newGuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewHlTarget, "", "", "e2e\$newGuid.md") | Out-Null
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"

$wsOverview.Range("G2").Value = $newHoGenerateDate

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $zhHandoffDate

# Latest Target File / Latest Handback File have not been produced yet for
# this new handoff, so clear them (and drop the old hyperlink on I2).
$wsZh.Hyperlinks.Delete()
$zhHlTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/$newGuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhHlTarget, "", "", "$newGuid.md") | Out-Null

$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $noHandbackDate

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newHoGenerateDate

$wsDe.Hyperlinks.Delete()
$deHlTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a2028cf7b0b4ea1100b4807308bddba0f46592e/e2e/$newGuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deHlTarget, "", "", "$newGuid.md") | Out-Null

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $noHandbackDate

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
